{"js": "// Replace the 25 division-problem prompts in the worksheet table.\n// The table has 5 \"data\" rows (0-based indices 0, 4, 8, 12, 16 \u2014 each\n// followed by 3 blank rows used for student work) with 5 columns\n// apiece. We walk the cells in document order (row-major) and swap the\n// old \"a\u00f7b=\" prompt for the new one, replacing just the run's text via\n// the paragraph range so the existing run/paragraph formatting\n// (TimeNewRoman, sz 30, left justification) is preserved untouched.\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"42\u00f76=\", newText: \"48\u00f79=\" },\n  { row: 0, col: 1, oldText: \"83\u00f72=\", newText: \"13\u00f73=\" },\n  { row: 0, col: 2, oldText: \"98\u00f75=\", newText: \"30\u00f74=\" },\n  { row: 0, col: 3, oldText: \"48\u00f75=\", newText: \"89\u00f78=\" },\n  { row: 0, col: 4, oldText: \"61\u00f77=\", newText: \"76\u00f79=\" },\n\n  { row: 4, col: 0, oldText: \"71\u00f79=\", newText: \"62\u00f73=\" },\n  { row: 4, col: 1, oldText: \"66\u00f76=\", newText: \"20\u00f75=\" },\n  { row: 4, col: 2, oldText: \"91\u00f72=\", newText: \"53\u00f72=\" },\n  { row: 4, col: 3, oldText: \"20\u00f77=\", newText: \"85\u00f72=\" },\n  { row: 4, col: 4, oldText: \"59\u00f73=\", newText: \"56\u00f76=\" },\n\n  { row: 8, col: 0, oldText: \"95\u00f74=\", newText: \"63\u00f79=\" },\n  { row: 8, col: 1, oldText: \"73\u00f76=\", newText: \"44\u00f76=\" },\n  { row: 8, col: 2, oldText: \"36\u00f77=\", newText: \"37\u00f75=\" },\n  { row: 8, col: 3, oldText: \"56\u00f78=\", newText: \"57\u00f79=\" },\n  { row: 8, col: 4, oldText: \"28\u00f77=\", newText: \"13\u00f79=\" },\n\n  { row: 12, col: 0, oldText: \"20\u00f77=\", newText: \"30\u00f75=\" },\n  { row: 12, col: 1, oldText: \"84\u00f76=\", newText: \"23\u00f78=\" },\n  { row: 12, col: 2, oldText: \"48\u00f74=\", newText: \"66\u00f74=\" },\n  { row: 12, col: 3, oldText: \"53\u00f76=\", newText: \"18\u00f77=\" },\n  { row: 12, col: 4, oldText: \"27\u00f74=\", newText: \"71\u00f79=\" },\n\n  { row: 16, col: 0, oldText: \"17\u00f77=\", newText: \"32\u00f74=\" },\n  { row: 16, col: 1, oldText: \"76\u00f79=\", newText: \"84\u00f75=\" },\n  { row: 16, col: 2, oldText: \"29\u00f73=\", newText: \"14\u00f75=\" },\n  { row: 16, col: 3, oldText: \"91\u00f75=\", newText: \"95\u00f72=\" },\n  { row: 16, col: 4, oldText: \"59\u00f75=\", newText: \"53\u00f74=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nfor (const rep of replacements) {\n  const current = table.values[rep.row][rep.col];\n  if (current !== rep.oldText) {\n    throw new Error(\n      `Unexpected text at row ${rep.row}, col ${rep.col}: expected \"${rep.oldText}\", found \"${current}\"`\n    );\n  }\n  const cell = table.getCell(rep.row, rep.col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(rep.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem prompts in the worksheet table.\n# The table has 5 \"data\" rows (each followed by 3 blank rows for student\n# work) with 5 columns apiece; we walk the cells in document order and\n# swap the old \"a\u00f7b=\" prompt for the new one. Matching against the\n# expected old value keeps this robust (and avoids accidentally touching\n# a cell whose value coincides with some other cell's replacement text).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{Row=1;  Col=1; Old=\"42\u00f76=\"; New=\"48\u00f79=\"},\n    @{Row=1;  Col=2; Old=\"83\u00f72=\"; New=\"13\u00f73=\"},\n    @{Row=1;  Col=3; Old=\"98\u00f75=\"; New=\"30\u00f74=\"},\n    @{Row=1;  Col=4; Old=\"48\u00f75=\"; New=\"89\u00f78=\"},\n    @{Row=1;  Col=5; Old=\"61\u00f77=\"; New=\"76\u00f79=\"},\n\n    @{Row=5;  Col=1; Old=\"71\u00f79=\"; New=\"62\u00f73=\"},\n    @{Row=5;  Col=2; Old=\"66\u00f76=\"; New=\"20\u00f75=\"},\n    @{Row=5;  Col=3; Old=\"91\u00f72=\"; New=\"53\u00f72=\"},\n    @{Row=5;  Col=4; Old=\"20\u00f77=\"; New=\"85\u00f72=\"},\n    @{Row=5;  Col=5; Old=\"59\u00f73=\"; New=\"56\u00f76=\"},\n\n    @{Row=9;  Col=1; Old=\"95\u00f74=\"; New=\"63\u00f79=\"},\n    @{Row=9;  Col=2; Old=\"73\u00f76=\"; New=\"44\u00f76=\"},\n    @{Row=9;  Col=3; Old=\"36\u00f77=\"; New=\"37\u00f75=\"},\n    @{Row=9;  Col=4; Old=\"56\u00f78=\"; New=\"57\u00f79=\"},\n    @{Row=9;  Col=5; Old=\"28\u00f77=\"; New=\"13\u00f79=\"},\n\n    @{Row=13; Col=1; Old=\"20\u00f77=\"; New=\"30\u00f75=\"},\n    @{Row=13; Col=2; Old=\"84\u00f76=\"; New=\"23\u00f78=\"},\n    @{Row=13; Col=3; Old=\"48\u00f74=\"; New=\"66\u00f74=\"},\n    @{Row=13; Col=4; Old=\"53\u00f76=\"; New=\"18\u00f77=\"},\n    @{Row=13; Col=5; Old=\"27\u00f74=\"; New=\"71\u00f79=\"},\n\n    @{Row=17; Col=1; Old=\"17\u00f77=\"; New=\"32\u00f74=\"},\n    @{Row=17; Col=2; Old=\"76\u00f79=\"; New=\"84\u00f75=\"},\n    @{Row=17; Col=3; Old=\"29\u00f73=\"; New=\"14\u00f75=\"},\n    @{Row=17; Col=4; Old=\"91\u00f75=\"; New=\"95\u00f72=\"},\n    @{Row=17; Col=5; Old=\"59\u00f75=\"; New=\"53\u00f74=\"}\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $rep.Old) {\n        throw \"Unexpected text at row $($rep.Row), col $($rep.Col): expected '$($rep.Old)', found '$current'\"\n    }\n    $cell.Range.Text = $rep.New\n}\n"}
